$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell (G1 - "sum") onto the new
# header cell H1, then set its text to "Save" - mirrors the style used by
# the other header cells.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Populate the new "Save" column values for rows 2-8.
$saveValues = @(1, 0, 0, 0, 0, 1, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}

$excel.CutCopyMode = 0
